# Edit script: Added Week 15 simulations
$wb = $excel.ActiveWorkbook

# --- YDS sheet: append week 15 simulation numbers to R/P distribution strings ---
$wsYDS = $wb.Worksheets.Item("YDS")
$wsYDS.Range("B2").Value = '1 3 1 2 17 2 7 -3 4 1 13 0 8 -1 6 5 2 1 0 -1 0 0 11 2 16 14 6 7 3 0 15 2 0 29 9 6 12 0 10 0 13 2 3 0 -2 8 1 1 9 7 1 3 1 4 8 2 2 12 2 2 9 3 4 3 6 3 -1 5 0 2 24 9 2 3 4 13 2 0 1 1 -2 2 3 1 11 2 4 1 0 1 10 5 -1 -2 2 1 6 6 -2 13 2 2 2 4 -2 -2 9 4 0 5 5 11 3 6 0 3 3 3 0 5 48 2 4 -1 3 2 2 1 2 1 0 -3 2 -3 3 4 14 5 9 1 3 4 3 5 1 8 9 6 0 14 2 0 30 16 -1 -1 23 4 0 3 3 -1 -3 16 5 4 16 6 2 -1 15 0 2 3 2 0 0 13 19 4 3 -2 8 4 8 8 7 1 4 -4 11 15 -3 4 66 1 3 3 0 0 1 2 -1 1 0 -1 0 -8 1 24 1 3 -1 5 3 -4 12 6 6 8 5 5 3 4 7 4 -3 5 -1 3 1 4 0 2 7 0 1 5 2 4 1 11 -1 6 -1 3 3 1 2 5 0 -1 4 0 7 4 3 5 13 2 4 0 7 2 0 4 3 7 12 4 1 3 1 1 0 4 7 4 30 0 6 -4 7 1 3 1 2 7 5 1 4 9 3 5 -4 3 10 -1 2 5 1 8 11 4 8 0 1 3 4 4 2 5 1 20 4 30 14 3 3 16 15 3 29 9 3 0 1 7 7 9 8 5 -4 7 14 -3 5 -2 8 5 4 1 11 2 4 -1 0 2'
$wsYDS.Range("B3").Value = '9 6 12 11 9 -1 2 9 4 17 11 25 5 5 6 20 0 9 10 10 5 6 9 8 34 10 7 27 0 6 14 7 3 5 12 64 6 14 5 3 7 9 14 -3 4 15 6 6 12 13 7 13 3 17 12 9 8 20 17 8 10 7 6 23 26 7 15 7 4 7 9 16 3 -4 15 5 8 9 -1 21 11 28 4 14 15 2 11 9 1 11 22 4 12 20 2 3 1 6 20 17 14 -3 31 5 3 19 6 37 14 4 17 6 3 6 14 22 2 15 4 9 2 11 6 3 14 9 7 16 8 21 6 19 11 9 6 7 7 3 3 20 7 24 13 2 -6 8 14 7 9 24 35 5 17 26 40 3 4 12 6 8 6 7 13 -3 27 31 3 20 3 5 7 -2 6 11 18 7 7 -4 -2 4 4 9 -3 6 13 6 3 6 50 5 3 7 7 20 4 18 5 9 5 6 27 1 2 14 6 -1 10 24 7 7 9 3 18 21 21 5 4 27 19 3 12 15 16 3 27 5 18 5 43 3 18 56 15 10 5 7 15 6 7 9 21 5 19 6 7 7 3 23 6 19 26 2 11 29 7 4 30 15 20 7 6 14 15 24 16 5 17 7 5 23 2 19 0 18 -2 -2 34 6 7 14 6 2 16 15 11 5 19 10 48 7 16 4 21 14 16 5 14 5 -2 11 3 8 5 8 14 18 22 19 17 13 0 1 62 12 17'
$wsYDS.Range("C2").Value = '5 -2 3 4 1 2 -1 13 1 3 6 4 7 2 2 5 10 19 1 5 5 9 12 -1 7 1 2 8 2 4 0 3 0 5 0 14 0 9 6 -1 6 2 8 12 8 9 3 4 6 3 1 6 6 -3 4 2 7 4 2 2 1 8 8 1 11 8 3 3 30 5 6 1 4 4 3 4 5 3 4 2 4 5 2 3 2 5 2 9 6 1 0 1 33 3 12 13 2 1 6 8 1 2 6 2 2 9 1 13 0 2 3 2 3 12 1 1 11 4 9 5 4 3 1 2 4 13 4 9 0 3 3 0 4 3 7 3 2 8 5 3 6 2 3 5 -2 6 9 6 0 5 2 5 2 6 30 0 13 8 2 6 2 3 5 4 2 7 2 5 2 2 2 6 2 4 2 2 4 2 1 2 1 5 8 5 6 4 4 2 6 3 4 9 3 5 6 7 11 6 12 0 1 7 1 4 9 4 6 1 3 8 7 8 5 10 1 6 5 4 2 3 19 5 12 1 5 4 5 1 0 4 3 17 1 1 3 5 5 9 2 6 3 3 5 4 1 5 3 2 1 0 11 2 5 4 18 11 7 6 5 7 4 4 3 -3 6 15 20 10 1 3 4 9 3 4 10 4 4 2 5 49 4 3 2 4 12 0 1 3 -6 5 1 0 6 7 9 4 3 1 2 10 1 2 8 4 5 4 6 8 2 4 3 0 2 -2 2 11 1 0 5 9 3 5 6 0 2 3 -2 7 2 -1 3 -2 9 3 4 4 23 8 3 11 13 1 1 2 3 6 -2'
$wsYDS.Range("C3").Value = '5 3 7 12 17 13 2 14 2 9 50 4 28 5 14 15 18 5 6 4 8 15 16 4 0 15 3 34 9 19 77 4 18 6 13 25 9 11 29 6 2 6 17 2 3 7 17 5 28 10 10 12 9 5 21 7 39 7 12 11 31 5 14 6 -3 19 7 17 15 3 11 21 8 3 3 9 8 9 12 11 5 9 24 19 10 8 13 4 13 8 4 4 6 21 8 14 5 2 17 2 8 7 18 21 10 9 0 16 18 17 5 15 14 1 4 4 21 6 4 7 5 8 17 8 20 29 4 9 3 7 73 10 2 35 1 13 15 7 19 33 18 1 7 -3 15 5 13 13 15 7 -1 9 4 5 12 22 3 3 8 8 18 11 5 22 6 8 7 7 20 20 7 13 7 11 10 18 2 7 8 5 4 22 17 11 10 2 7 7 26 13 6 2 37 2 4 13 3 5 39 8 17 12 15 8 25 13 12 17 10 9 13 26 4 18 75 21 37 3 2 11 2 9 24 3 13 13 30 12 8 11 5 26 9 9 25 9 28 23 6 5 27 7 6 24 7 8 14 9 5 6 13 9 9 9 10 8 7 6 15 6 5 6 6 17 29 5 -2 8 3 14 1 37 30 6 15 6 38 9 1 13 3 10 12'

# --- OFF sheet: updated totals ---
$wsOFF = $wb.Worksheets.Item("OFF")
$wsOFF.Range("C2").Value = 175
$wsOFF.Range("D2").Value = 9
$wsOFF.Range("F2").Value = 56
$wsOFF.Range("G2").Value = 47
$wsOFF.Range("I2").Value = 2
$wsOFF.Range("J2").Value = 21
$wsOFF.Range("L2").Value = 213
$wsOFF.Range("M2").Value = 136
$wsOFF.Range("O2").Value = 17
$wsOFF.Range("P2").Value = 10
$wsOFF.Range("Q2").Value = 392
$wsOFF.Range("B3").Value = 13
$wsOFF.Range("C3").Value = 129
$wsOFF.Range("E3").Value = 38
$wsOFF.Range("F3").Value = 69
$wsOFF.Range("G3").Value = 27
$wsOFF.Range("H3").Value = 32
$wsOFF.Range("I3").Value = 51
$wsOFF.Range("J3").Value = 42

# --- DEF sheet: updated totals ---
$wsDEF = $wb.Worksheets.Item("DEF")
$wsDEF.Range("C2").Value = 183
$wsDEF.Range("D2").Value = 7
$wsDEF.Range("E2").Value = 9
$wsDEF.Range("F2").Value = 50
$wsDEF.Range("G2").Value = 45
$wsDEF.Range("J2").Value = 18
$wsDEF.Range("L2").Value = 214
$wsDEF.Range("M2").Value = 136
$wsDEF.Range("O2").Value = 17
$wsDEF.Range("P2").Value = 11
$wsDEF.Range("Q2").Value = 387
$wsDEF.Range("B3").Value = 9
$wsDEF.Range("C3").Value = 135
$wsDEF.Range("E3").Value = 19
$wsDEF.Range("F3").Value = 96
$wsDEF.Range("G3").Value = 36
$wsDEF.Range("H3").Value = 26
$wsDEF.Range("I3").Value = 52
$wsDEF.Range("J3").Value = 37
$wsDEF.Range("N3").Value = 19

# --- ST sheet: updated totals + appended week 15 simulation strings ---
$wsST = $wb.Worksheets.Item("ST")
$wsST.Range("B2").Value = 77
$wsST.Range("D2").Value = 54
$wsST.Range("F2").Value = 69
$wsST.Range("G2").Value = 61
$wsST.Range("J2").Value = 65
$wsST.Range("K2").Value = 54
$wsST.Range("L2").Value = 35
$wsST.Range("M2").Value = 22
$wsST.Range("N2").Value = 14
$wsST.Range("B3").Value = 64
$wsST.Range("B4").Value = '66 45 62 50 65 69 67 61 67 64 69 65 67'
$wsST.Range("B5").Value = '11 14 21 21 26 24 22 24 14 30 23 33 17'
$wsST.Range("D3").Value = '39 46 49 61 52 41 51 63 46 43 50 60 32 49 45 58 47 58 46 36 46 50 54 46 44 54 53 37 44 37 37 39 51 29 53 54 41 53 41 40 42 55 45 33 47 45 56 59 46 49 37 46 49 51'
$wsST.Range("D4").Value = '8 15 0 12 15 3 0 0 6 0 8 5 0 0 0 8 8 9 0 2 -2 14 0 0 1 0 4 0 0 7 0 -3 0 0 0 0 0 14 0 9 14 13 0 0 0 0 12 2 0 0 0 0 8 0'
$wsST.Range("D5").Value = '0 0 0 0 0 0 12 18 0 -4 15 13 -1 0 11 0 0 2 0 0 0 0 2 0 7 0 0 0 0 0 7 0 0 0 0 0 45 0 0 7 17 0 8 0 0 0 16 2 5 16 12 0 0 15 0'

# --- TURNS sheet: updated totals ---
$wsTURNS = $wb.Worksheets.Item("TURNS")
$wsTURNS.Range("B2").Value = 6
$wsTURNS.Range("C2").Value = 3
$wsTURNS.Range("E2").Value = 6
$wsTURNS.Range("D3").Value = 4
$wsTURNS.Range("E3").Value = 5

# --- PEN sheet: updated totals ---
$wsPEN = $wb.Worksheets.Item("PEN")
$wsPEN.Range("D2").Value = 7
$wsPEN.Range("D3").Value = 6
$wsPEN.Range("D4").Value = 8

Write-Output "Applied Week 15 simulations update"
